$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,('spring', '春|はる')
    ,('spring roll', '春巻|はるまき')
    ,('vernal equinox', '春分|しゅんぶん')
    ,('youth', '青春|せいしゅん')
    ,('autumn', '秋|あき')
    ,('fall semester', '秋学期|あきがっき')
    ,('autumnal equinox', '秋分|しゅうぶん')
    ,('winter', '冬|ふゆ')
    ,('winter vacation', '冬休み|ふゆやすみ')
    ,('warm winter', '暖冬|だんとう')
    ,('four seasons', '春夏秋冬|しゅんかしゅうとう')
    ,('flower', '花|はな')
    ,('flower viewing', '花見|はなみ')
    ,('fireworks', '花火|はなび')
    ,('hay fever', '花粉症|かふんしょう')
    ,('flower arrangement', '生け花|いけばな')
    ,('Mr./Ms.', '～様|～さま')
    ,('customer (respectful)', 'お客様|おきゃくさま')
    ,('everyone', '皆様|みなさま')
    ,('manner', '様子|ようす')
    ,('uneasy; worried', '不安な|ふあんな')
    ,('recession', '不景気|ふけいき')
    ,('inconvenience', '不便|ふべん')
    ,('lack of water', '水不足|みずぶそく')
    ,('older sister', '姉|あね')
    ,('older sister (honorific)', 'お姉さん|おねえさん')
    ,('sisters', '姉妹|しまい')
    ,('sister city', '姉妹都市|しまいとし')
    ,('older brother', '兄|あに')
    ,('older brother (honorific)', 'お兄さん|おにいさん')
    ,('brothers', '兄弟|きょうだい')
    ,('Chinese character', '漢字|かんじ')
    ,('Chinese herbal medicine', '漢方薬|かんぽうやく')
    ,('Kanji dictionary', '漢和辞典|かんわじてん')
    ,('to graduate', '卒業する|そつぎょうする')
    ,('graduation ceremony', '卒業式|そつぎょうしき')
    ,('university graduate', '大卒|だいそつ')
    ,('engineering', '工学|こうがく')
    ,('construction', '工事|こうじ')
    ,('factory', '工場|こうじょう')
    ,('carpenter', '大工|だいく')
    ,('research', '研究|けんきゅう')
    ,('researcher', '研究者|けんきゅうしゃ')
    ,('training', '研修|けんしゅう')
    ,('to sharpen', '研ぐ|とぐ')
    ,('inquiry', '探究|たんきゅう')
    ,('to investigate thoroughly', '究める|きわめる')
    ,('question', '質問|しつもん')
    ,('good quality', '質がいい|しつがいい')
    ,('aptitude', '素質|そしつ')
    ,('pawn shop', '質屋|しちや')
    ,('problem', '問題|もんだい')
    ,('visit', '訪問|ほうもん')
    ,('to question', '問う|とう')
    ,('many', '多い|おおい')
    ,('probably', '多分|たぶん')
    ,('many of...', '～の多く|～のおおく')
    ,('decision by majority', '多数決|たすうけつ')
)

$startRow = 53
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
